# Update header labels and repeated labels to their new (English) values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "id"
$ws.Range("B1").Value = "record"
$ws.Range("C1").Value = "address"

$ws.Range("A6").Value = "id"
$ws.Range("A7").Value = "address"
$ws.Range("A8").Value = "record"

# Move the active selection to A8, matching the saved workbook state.
$ws.Range("A8").Select()
